$wb = $excel.ActiveWorkbook

# --- Populate Sheet2 with the login credentials used by the Sikuli test cases ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "EmailAddress"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "swbtop@gmail.com"
$ws2.Range("B2").Value = "admin"

# The header/data cells use the same text-format style (numFmtId 49) as Sheet1.
$ws2.Range("A1:B2").NumberFormat = "@"

# Auto-fit column A like Sheet1's columns were auto-fit to their content.
$ws2.Columns.Item(1).AutoFit()

# Sheet2's table should print in portrait orientation.
$ws2.PageSetup.Orientation = 1

# Sheet2 becomes the active / selected sheet with A1:B2 selected.
$ws2.Activate()
$ws2.Range("A1:B2").Select()

$wb.Save()
